$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Load Points")

$ws.Range("F2").Value = 0.5786482193441778
$ws.Range("G2").Value = 134
$ws.Range("H2").Value = 2.59096217616796
$ws.Range("I2").Value = 0.2233333333333333
$ws.Range("J2").Value = 185.6570502286639
$ws.Range("K2").Value = 46.9
$ws.Range("L2").Value = 121.5161260622773
$ws.Range("M2").Value = 544.1020569952716
$ws.Range("N2").Value = 0.3094284170477732
$ws.Range("F3").Value = 0.6580166761854783
$ws.Range("G3").Value = 142
$ws.Range("H3").Value = 2.780352152896387
$ws.Range("I3").Value = 0.2366666666666667
$ws.Range("J3").Value = 196.6509103318847
$ws.Range("K3").Value = 49.7
$ws.Range("L3").Value = 138.1835019989504
$ws.Range("M3").Value = 583.8739521082413
$ws.Range("N3").Value = 0.3277515172198078
$ws.Range("F4").Value = 0.9831601997505144
$ws.Range("G4").Value = 148
$ws.Range("H4").Value = 3.985784593583166
$ws.Range("I4").Value = 0.2466666666666667
$ws.Range("J4").Value = 287.3517655477486
$ws.Range("K4").Value = 51.8
$ws.Range("L4").Value = 206.463641947608
$ws.Range("M4").Value = 837.0147646524649
$ws.Range("N4").Value = 0.4789196092462477
$ws.Range("F5").Value = 0.6512537367033345
$ws.Range("G5").Value = 136
$ws.Range("H5").Value = 2.87317825016177
$ws.Range("I5").Value = 0.2266666666666667
$ws.Range("J5").Value = 238.334655591122
$ws.Range("K5").Value = 0.2266666666666667
$ws.Range("L5").Value = 0.6512537367033345
$ws.Range("M5").Value = 2.87317825016177
$ws.Range("N5").Value = 0.3972244259852034
$ws.Range("F6").Value = 0.7736264829465828
$ws.Range("G6").Value = 148
$ws.Range("H6").Value = 3.136323579513173
$ws.Range("I6").Value = 0.2466666666666667
$ws.Range("J6").Value = 230.3646322165773
$ws.Range("K6").Value = 0.2466666666666667
$ws.Range("L6").Value = 0.7736264829465828
$ws.Range("M6").Value = 3.136323579513173
$ws.Range("N6").Value = 0.3839410536942955
$ws.Range("F7").Value = 0.6571696678955942
$ws.Range("G7").Value = 144
$ws.Range("H7").Value = 2.738206949564976
$ws.Range("I7").Value = 0.24
$ws.Range("J7").Value = 200.8457690670942
$ws.Range("K7").Value = 2.4
$ws.Range("L7").Value = 6.571696678955941
$ws.Range("M7").Value = 27.38206949564976
$ws.Range("N7").Value = 0.334742948445157
$ws.Range("F8").Value = 0.6971688234002389
$ws.Range("G8").Value = 148
$ws.Range("H8").Value = 2.826360094865833
$ws.Range("I8").Value = 0.2466666666666667
$ws.Range("J8").Value = 135.8568492431152
$ws.Range("K8").Value = 2.466666666666667
$ws.Range("L8").Value = 6.971688234002389
$ws.Range("M8").Value = 28.26360094865833
$ws.Range("N8").Value = 0.2264280820718586
$ws.Range("F9").Value = 0.6236740100605981
$ws.Range("G9").Value = 113
$ws.Range("H9").Value = 3.311543416250963
$ws.Range("I9").Value = 0.1883333333333333
$ws.Range("J9").Value = 178.3010857224789
$ws.Range("K9").Value = 0.1883333333333333
$ws.Range("L9").Value = 0.6236740100605981
$ws.Range("M9").Value = 3.311543416250963
$ws.Range("N9").Value = 0.2971684762041314
$ws.Range("F10").Value = 0.7831049390346055
$ws.Range("G10").Value = 113
$ws.Range("H10").Value = 4.15807932230764
$ws.Range("I10").Value = 0.1883333333333333
$ws.Range("J10").Value = 266.6544609368764
$ws.Range("K10").Value = 0.1883333333333333
$ws.Range("L10").Value = 0.7831049390346055
$ws.Range("M10").Value = 4.15807932230764
$ws.Range("N10").Value = 0.4444241015614606
$ws.Range("F11").Value = 0.8083704059571513
$ws.Range("G11").Value = 162
$ws.Range("H11").Value = 2.993964466507967
$ws.Range("I11").Value = 0.27
$ws.Range("J11").Value = 172.0179200339518
$ws.Range("K11").Value = 56.7
$ws.Range("L11").Value = 169.7577852510018
$ws.Range("M11").Value = 628.7325379666731
$ws.Range("N11").Value = 0.2866965333899197
$ws.Range("F12").Value = 0.7834668598912301
$ws.Range("G12").Value = 159
$ws.Range("H12").Value = 2.956478716570679
$ws.Range("I12").Value = 0.265
$ws.Range("J12").Value = 246.9406585675197
$ws.Range("K12").Value = 55.65000000000001
$ws.Range("L12").Value = 164.5280405771583
$ws.Range("M12").Value = 620.8605304798426
$ws.Range("N12").Value = 0.4115677642791996
$ws.Range("F13").Value = 0.8298225553951707
$ws.Range("G13").Value = 160
$ws.Range("H13").Value = 3.11183458273189
$ws.Range("I13").Value = 0.2666666666666667
$ws.Range("J13").Value = 254.0287825301938
$ws.Range("K13").Value = 53.33333333333334
$ws.Range("L13").Value = 165.9645110790341
$ws.Range("M13").Value = 622.3669165463781
$ws.Range("N13").Value = 0.4233813042169896
$ws.Range("F14").Value = 0.764627860346198
$ws.Range("G14").Value = 170
$ws.Range("H14").Value = 2.698686565927758
$ws.Range("I14").Value = 0.2833333333333333
$ws.Range("J14").Value = 227.2056409479531
$ws.Range("K14").Value = 0.2833333333333333
$ws.Range("L14").Value = 0.764627860346198
$ws.Range("M14").Value = 2.698686565927758
$ws.Range("N14").Value = 0.3786760682465886
$ws.Range("F15").Value = 0.661600576512557
$ws.Range("G15").Value = 173
$ws.Range("H15").Value = 2.294568473453955
$ws.Range("I15").Value = 0.2883333333333333
$ws.Range("J15").Value = 210.2484154518459
$ws.Range("K15").Value = 0.2883333333333333
$ws.Range("L15").Value = 0.661600576512557
$ws.Range("M15").Value = 2.294568473453955
$ws.Range("N15").Value = 0.3504140257530765
$ws.Range("F16").Value = 0.7935588273947127
$ws.Range("G16").Value = 161
$ws.Range("H16").Value = 2.957362089669737
$ws.Range("I16").Value = 0.2683333333333333
$ws.Range("J16").Value = 157.7231527787888
$ws.Range("K16").Value = 2.683333333333333
$ws.Range("L16").Value = 7.935588273947127
$ws.Range("M16").Value = 29.57362089669737
$ws.Range("N16").Value = 0.2628719212979814
$ws.Range("F17").Value = 0.7328048259383743
$ws.Range("G17").Value = 138
$ws.Range("H17").Value = 3.186107938862496
$ws.Range("I17").Value = 0.23
$ws.Range("J17").Value = 170.8188476392567
$ws.Range("K17").Value = 2.3
$ws.Range("L17").Value = 7.328048259383743
$ws.Range("M17").Value = 31.86107938862497
$ws.Range("N17").Value = 0.2846980793987611
$ws.Range("F18").Value = 0.6891439518381779
$ws.Range("G18").Value = 142
$ws.Range("H18").Value = 2.911875852837372
$ws.Range("I18").Value = 0.2366666666666667
$ws.Range("J18").Value = 163.2701446236129
$ws.Range("K18").Value = 47.33333333333334
$ws.Range("L18").Value = 137.8287903676356
$ws.Range("M18").Value = 582.3751705674744
$ws.Range("N18").Value = 0.2721169077060215
$ws.Range("F19").Value = 0.6058171214429496
$ws.Range("G19").Value = 134
$ws.Range("H19").Value = 2.712613976610222
$ws.Range("I19").Value = 0.2233333333333333
$ws.Range("J19").Value = 174.3352134469874
$ws.Range("K19").Value = 44.66666666666666
$ws.Range("L19").Value = 121.1634242885899
$ws.Range("M19").Value = 542.5227953220444
$ws.Range("N19").Value = 0.2905586890783123
$ws.Range("F20").Value = 0.6144837694937428
$ws.Range("G20").Value = 138
$ws.Range("H20").Value = 2.671668563016273
$ws.Range("I20").Value = 0.23
$ws.Range("J20").Value = 169.6146098924079
$ws.Range("K20").Value = 46
$ws.Range("L20").Value = 122.8967538987486
$ws.Range("M20").Value = 534.3337126032546
$ws.Range("N20").Value = 0.2826910164873466
$ws.Range("F21").Value = 0.6664480209996395
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = 2.961991204442842
$ws.Range("I21").Value = 0.225
$ws.Range("J21").Value = 174.9747183350294
$ws.Range("K21").Value = 0.225
$ws.Range("L21").Value = 0.6664480209996395
$ws.Range("M21").Value = 2.961991204442842
$ws.Range("N21").Value = 0.2916245305583824
$ws.Range("F22").Value = 0.6906647426133448
$ws.Range("G22").Value = 148
$ws.Range("H22").Value = 2.79999219978383
$ws.Range("I22").Value = 0.2466666666666667
$ws.Range("J22").Value = 177.7831288087303
$ws.Range("K22").Value = 0.2466666666666667
$ws.Range("L22").Value = 0.6906647426133448
$ws.Range("M22").Value = 2.79999219978383
$ws.Range("N22").Value = 0.2963052146812172
$ws.Range("F23").Value = 0.5958149703453176
$ws.Range("G23").Value = 141
$ws.Range("H23").Value = 2.535382852533266
$ws.Range("I23").Value = 0.235
$ws.Range("J23").Value = 163.0084991053756
$ws.Range("K23").Value = 2.35
$ws.Range("L23").Value = 5.958149703453176
$ws.Range("M23").Value = 25.35382852533267
$ws.Range("N23").Value = 0.2716808318422927
$ws.Range("K24").Value = 0.2443273934311671
$ws.Range("L24").Value = 0.7278211462211549
$ws.Range("M24").Value = 2.978876563942061
$ws.Range("N24").Value = 7.303311518412024
$ws.Range("Q24").Value = 0.04216931777475056
